# Fix off-by-one bug in the symbol table / memory map address computation:
# the "end address" column (column B) for several memory segments was showing
# an address ending in ...FFC instead of the correct ...FFF value.
# This corresponds to the commit's "fixed bug when computing addresses of the
# variables for the symbol table".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value  = "0xFFFFFFFF"
$ws.Range("B6").Value  = "0xFFFFF7FF"
$ws.Range("B8").Value  = "0xFFFFEFFF"
$ws.Range("B10").Value = "0xBFFFFFFF"
$ws.Range("B12").Value = "0x10000FFF"
$ws.Range("B14").Value = "0x0FFFFFFF"
$ws.Range("B16").Value = "0x0000FFFF"
